# Auto-generated edit script applying numeric updates described in the commit diff
# for Sheets/Brynhildr_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 98.045456
$ws.Range("I9").Value = 98.111115
$ws.Range("K9").Value = 98.111115
$ws.Range("M9").Value = 70.888885

$ws.Range("H32").Value = 1933.1428
$ws.Range("I32").Value = 1928.2
$ws.Range("J32").Value = 1945.5
$ws.Range("K32").Value = 1928.2
$ws.Range("L32").Value = 1945.5
$ws.Range("M32").Value = -1602.2
$ws.Range("N32").Value = -2597.5

$ws.Range("H76").Value = 3943.2258
$ws.Range("I76").Value = 3852.2778
$ws.Range("J76").Value = 4069.1538
$ws.Range("K76").Value = 3852.2778
$ws.Range("L76").Value = 4069.1538
$ws.Range("M76").Value = -3537.2778
$ws.Range("N76").Value = -4699.1538

$ws.Range("H79").Value = 3943.2258
$ws.Range("I79").Value = 3852.2778
$ws.Range("J79").Value = 4069.1538
$ws.Range("K79").Value = 3852.2778
$ws.Range("L79").Value = 4069.1538
$ws.Range("M79").Value = -2760.2778
$ws.Range("N79").Value = -6253.1538

$ws.Range("H88").Value = 2253.72
$ws.Range("I88").Value = 3881.7778
$ws.Range("J88").Value = 1337.9375
$ws.Range("K88").Value = 3881.7778
$ws.Range("L88").Value = 1337.9375
$ws.Range("M88").Value = -3475.7778
$ws.Range("N88").Value = -2149.9375

$ws.Range("H91").Value = 2253.72
$ws.Range("I91").Value = 3881.7778
$ws.Range("J91").Value = 1337.9375
$ws.Range("K91").Value = 3881.7778
$ws.Range("L91").Value = 1337.9375
$ws.Range("M91").Value = -2477.7778
$ws.Range("N91").Value = -4145.9375

$ws.Range("H107").Value = 2011.3793
$ws.Range("I107").Value = 2228.6155
$ws.Range("K107").Value = 2228.6155
$ws.Range("M107").Value = -308.6154999999999

$ws.Range("H111").Value = 592.5714
$ws.Range("I111").Value = 562.25
$ws.Range("J111").Value = 633
$ws.Range("K111").Value = 1686.75
$ws.Range("L111").Value = 1899
$ws.Range("M111").Value = 1380.25
$ws.Range("N111").Value = -8033

$ws.Range("H112").Value = 1970.4286
$ws.Range("I112").Value = 2073.75
$ws.Range("J112").Value = 1832.6666
$ws.Range("K112").Value = 6221.25
$ws.Range("L112").Value = 5497.9998
$ws.Range("M112").Value = -5113.25
$ws.Range("N112").Value = -7713.9998

$ws.Range("H117").Value = 70707
$ws.Range("J117").Value = 70707
$ws.Range("L117").Value = 70707
$ws.Range("N117").Value = -79885

$ws.Range("H138").Value = 4642.222
$ws.Range("I138").Value = 4254.2856
$ws.Range("K138").Value = 12762.8568
$ws.Range("M138").Value = -7622.856800000001

$ws.Range("H141").Value = 11223.425
$ws.Range("I141").Value = 4209.7856
$ws.Range("K141").Value = 12629.3568
$ws.Range("M141").Value = -7449.356800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 7064.7
$ws.Range("I132").Value = 5566.6665
$ws.Range("J132").Value = 7706.7144
$ws.Range("K132").Value = 16699.9995
$ws.Range("L132").Value = 23120.1432
$ws.Range("M132").Value = -14169.9995
$ws.Range("N132").Value = -28180.1432

$ws.Range("H137").Value = 68231.664
$ws.Range("J137").Value = 68231.664
$ws.Range("L137").Value = 68231.664
$ws.Range("N137").Value = -78431.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2386.853
$ws.Range("I86").Value = 2247
$ws.Range("J86").Value = 2643.25
$ws.Range("K86").Value = 2247
$ws.Range("L86").Value = 2643.25
$ws.Range("M86").Value = -1124
$ws.Range("N86").Value = -4889.25

$ws.Range("H89").Value = 2386.853
$ws.Range("I89").Value = 2247
$ws.Range("J89").Value = 2643.25
$ws.Range("K89").Value = 11235
$ws.Range("L89").Value = 13216.25
$ws.Range("M89").Value = -5619
$ws.Range("N89").Value = -24448.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 29800
$ws.Range("J50").Value = 29800
$ws.Range("L50").Value = 29800
$ws.Range("N50").Value = -31050

$ws.Range("H51").Value = 29500
$ws.Range("J51").Value = 29500
$ws.Range("L51").Value = 29500
$ws.Range("N51").Value = -30972

$ws.Range("H60").Value = 14300
$ws.Range("I60").Value = 9000
$ws.Range("J60").Value = 15625
$ws.Range("K60").Value = 9000
$ws.Range("L60").Value = 15625
$ws.Range("M60").Value = -8489
$ws.Range("N60").Value = -16647

$ws.Range("H61").Value = 29500
$ws.Range("J61").Value = 29500
$ws.Range("L61").Value = 29500
$ws.Range("N61").Value = -30196

$ws.Range("H134").Value = 3084.8333
$ws.Range("I134").Value = 3062.2424
$ws.Range("J134").Value = 3333.3333
$ws.Range("K134").Value = 9186.727200000001
$ws.Range("L134").Value = 9999.999899999999
$ws.Range("M134").Value = -6651.727200000001
$ws.Range("N134").Value = -15069.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1966136.9
$ws.Range("J5").Value = 3876830
$ws.Range("L5").Value = 11630490
$ws.Range("N5").Value = -11630714

$ws.Range("H15").Value = 333396.66
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H68").Value = 6619.2905
$ws.Range("I68").Value = 5349.75
$ws.Range("J68").Value = 6807.3706
$ws.Range("K68").Value = 16049.25
$ws.Range("L68").Value = 20422.1118
$ws.Range("M68").Value = -15238.25
$ws.Range("N68").Value = -22044.1118

$ws.Range("H71").Value = 6619.2905
$ws.Range("I71").Value = 5349.75
$ws.Range("J71").Value = 6807.3706
$ws.Range("K71").Value = 48147.75
$ws.Range("L71").Value = 61266.3354
$ws.Range("M71").Value = -44091.75
$ws.Range("N71").Value = -69378.33540000001

$ws.Range("H86").Value = 707.5

$ws.Range("H89").Value = 707.5

$ws.Range("H117").Value = 8924.091
$ws.Range("I117").Value = 157.5
$ws.Range("J117").Value = 10872.223
$ws.Range("K117").Value = 472.5
$ws.Range("L117").Value = 32616.669
$ws.Range("M117").Value = 2969.5
$ws.Range("N117").Value = -39500.669

$ws.Range("H126").Value = 12142.857
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -10060

$ws.Range("H135").Value = 1966136.9
$ws.Range("J135").Value = 3876830
$ws.Range("L135").Value = 34891470
$ws.Range("N135").Value = -34896540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2046
$ws.Range("I102").Value = 840.3684
$ws.Range("J102").Value = 9681.666999999999
$ws.Range("K102").Value = 840.3684
$ws.Range("L102").Value = 9681.666999999999
$ws.Range("M102").Value = 781.6316
$ws.Range("N102").Value = -12925.667

$ws.Range("H107").Value = 1212.5
$ws.Range("I107").Value = 1416.6666
$ws.Range("K107").Value = 1416.6666
$ws.Range("M107").Value = 503.3334

$ws.Range("H113").Value = 1499.9166
$ws.Range("I113").Value = 1499.9166
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1499.9166
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 670.0834
$ws.Range("N113").ClearContents()

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 109500
$ws.Range("J36").Value = 109500
$ws.Range("L36").Value = 109500
$ws.Range("N36").Value = -110624

$ws.Range("H40").Value = 4169.524
$ws.Range("I40").Value = 2897.8462
$ws.Range("K40").Value = 2897.8462
$ws.Range("M40").Value = -2761.8462

$ws.Range("H55").Value = 1624.3448
$ws.Range("I55").Value = 1674.375
$ws.Range("J55").Value = 1605.2858
$ws.Range("K55").Value = 1674.375
$ws.Range("L55").Value = 1605.2858
$ws.Range("M55").Value = -1501.375
$ws.Range("N55").Value = -1951.2858

$ws.Range("H61").Value = 10900.454
$ws.Range("I61").Value = 13685.714
$ws.Range("K61").Value = 13685.714
$ws.Range("M61").Value = -13483.714

$ws.Range("H82").Value = 102000
$ws.Range("I82").Value = 102000
$ws.Range("K82").Value = 102000
$ws.Range("M82").Value = -101639

$ws.Range("H85").Value = 102000
$ws.Range("I85").Value = 102000
$ws.Range("K85").Value = 102000
$ws.Range("M85").Value = -100752

$ws.Range("H113").Value = 10900.454
$ws.Range("I113").Value = 13685.714
$ws.Range("K113").Value = 13685.714
$ws.Range("M113").Value = -11515.714

$ws.Range("H122").Value = 5057
$ws.Range("I122").Value = 3863.5
$ws.Range("J122").Value = 7444
$ws.Range("K122").Value = 11590.5
$ws.Range("L122").Value = 22332
$ws.Range("M122").Value = -9140.5
$ws.Range("N122").Value = -27232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 52500
$ws.Range("J106").Value = 52500
$ws.Range("L106").Value = 52500
$ws.Range("N106").Value = -55024

$ws.Range("H122").Value = 44758.69
$ws.Range("I122").Value = 1189.7368
$ws.Range("J122").Value = 163017.28
$ws.Range("K122").Value = 3569.2104
$ws.Range("L122").Value = 489051.84
$ws.Range("M122").Value = -1119.2104
$ws.Range("N122").Value = -493951.84

$ws.Range("H126").Value = 1437.6923
$ws.Range("I126").Value = 1163.6364
$ws.Range("J126").Value = 2945
$ws.Range("K126").Value = 3490.9092
$ws.Range("L126").Value = 8835
$ws.Range("M126").Value = -1020.9092
$ws.Range("N126").Value = -13775

$ws.Range("H132").Value = 9806367
$ws.Range("I132").Value = 12822681
$ws.Range("K132").Value = 38468043
$ws.Range("M132").Value = -38465513
